$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 0.00257494942063638
$ws.Range("J3").Value = 0.007540923303292256
$ws.Range("G4").Value = 0.02299061982711054
$ws.Range("H4").Value = 0.0006437373551590951
$ws.Range("G5").Value = 0.0003678499172337686
$ws.Range("H5").Value = 0.001379437189626632
$ws.Range("I5").Value = 0.03577340445098399
$ws.Range("J5").Value = 0.001655324627551959
$ws.Range("E6").Value = 0.08129483170866286
$ws.Range("F6").Value = 0.02151922015817546
$ws.Range("G6").Value = 0.2315615228986573
$ws.Range("H6").Value = 0.4261541291153209
$ws.Range("I6").Value = 0.0001839249586168843
$ws.Range("J6").Value = 0.01747287106860401
$ws.Range("E7").Value = 0.01747287106860401
$ws.Range("F7").Value = 0.04294647783704249
$ws.Range("G7").Value = 0.04193489056464962
$ws.Range("H7").Value = 0.001563362148243517
$ws.Range("I7").Value = 0.002299061982711054
$ws.Range("J7").Value = 0.0003678499172337686
$ws.Range("F8").Value = 0
$ws.Range("H8").Value = 0.005425786279198087
$ws.Range("E9").Value = 0.007081110906750046
$ws.Range("F9").Value = 0.1370240941695788
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0.009380172889461099
$ws.Range("J10").Value = 0.007173073386058488
$ws.Range("H11").Value = 0.01563362148243517
$ws.Range("J11").Value = 0.008460548096376677
$ws.Range("E12").Value = 0.0009196247930844216
$ws.Range("F12").Value = 0.0007356998344675372
$ws.Range("G12").Value = 0.001103549751701306
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("E13").Value = 0.0006437373551590951
$ws.Range("F13").Value = 0.005885598675740298
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0.000551774875850653
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0.003862424130954571
$ws.Range("E16").Value = 0.001747287106860401
$ws.Range("F16").Value = 0.002299061982711054
$ws.Range("G16").Value = 0.001655324627551959
$ws.Range("H16").Value = 0.0001839249586168843
$ws.Range("I16").Value = 0.001103549751701306
$ws.Range("J16").Value = 0.0007356998344675372
$ws.Range("E17").Value = 0.6103549751701305
$ws.Range("F17").Value = 0.4775611550487401
$ws.Range("G17").Value = 0.3601250689718595
$ws.Range("H17").Value = 0.1102630126908221
$ws.Range("I17").Value = 0.3673901048372264
$ws.Range("J17").Value = 0.4462019496045613
$ws.Range("E18").Value = 0.0004598123965422108
$ws.Range("F18").Value = 0.0004598123965422107
$ws.Range("G18").Value = 0.001655324627551959
$ws.Range("H18").Value = 0.000551774875850653
$ws.Range("I18").Value = 0.001655324627551959
$ws.Range("J18").Value = 0.0009196247930844216
$ws.Range("J19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0.0001839249586168843
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0.0003678499172337686
$ws.Range("E21").Value = 0.08295015633621483
$ws.Range("F21").Value = 0.1484274416038256
$ws.Range("G21").Value = 0.1427257678867022
$ws.Range("H21").Value = 0.1910060695236344
$ws.Range("I21").Value = 0.383943351112746
$ws.Range("J21").Value = 0.09085892955674085
$ws.Range("J22").Value = 0.03209490527864631
$ws.Range("J23").Value = 0.004414199006805224
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0.01167923487217215
$ws.Range("H24").Value = 0.07540923303292256
$ws.Range("J24").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0.001839249586168843
$ws.Range("I25").Value = 0.01232297222733125
$ws.Range("J25").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0.05949972411256207
$ws.Range("H26").Value = 0.006989148427441604
$ws.Range("I26").Value = 0.04782048924038992
$ws.Range("J26").Value = 0.1006069523634357
$ws.Range("E27").Value = 0.8942431487952915
$ws.Range("F27").Value = 0.7611734412359757
$ws.Range("G27").Value = 0.1419900680522347
$ws.Range("H27").Value = 0.2223652749678131
$ws.Range("I27").Value = 0.5663049475813867
$ws.Range("J27").Value = 0.4908037520691558
$ws.Range("E28").Value = 0.08754828030163693
$ws.Range("F28").Value = 0.232113297774508
$ws.Range("G28").Value = 0.01305867206179879
$ws.Range("H28").Value = 0.1392311936729814
$ws.Range("I28").Value = 0.1565201397829686
$ws.Range("J28").Value = 0.3133161670038624
$ws.Range("E29").Value = 0.0007356998344675372
$ws.Range("F29").Value = 0.0009196247930844216
$ws.Range("G29").Value = 0.6529336030899393
$ws.Range("H29").Value = 0.5442339525473607
$ws.Range("I29").Value = 0.001471399668935074
$ws.Range("J29").Value = 0.000551774875850653
$ws.Range("E30").Value = 0.0003678499172337686
$ws.Range("F30").Value = 0.0007356998344675372
$ws.Range("G30").Value = 0.0002758874379253265
$ws.Range("H30").Value = 0.007173073386058488
$ws.Range("I30").Value = 0.00386242413095457
$ws.Range("J30").Value = 0.0003678499172337686
$ws.Range("E31").Value = 0.000551774875850653
$ws.Range("F31").Value = 0.0004598123965422108
$ws.Range("G31").Value = 0.000551774875850653
$ws.Range("H31").Value = 0.01066764759977929
$ws.Range("I31").Value = 0.006253448592974067
$ws.Range("J31").Value = 0.0006437373551590951
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0.02547360676843848
$ws.Range("G33").Value = 0.03512966709582491
$ws.Range("H33").Value = 0.04073937833363987
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("E34").Value = 0.01121942247562994
$ws.Range("F34").Value = 0.002023174544785727
$ws.Range("J34").Value = 0
$ws.Range("G35").Value = 0.0003678499172337686
$ws.Range("H35").Value = 0.002666911899944822
$ws.Range("I35").Value = 0.01278278462387346
$ws.Range("J35").Value = 0.03724480411991907
$ws.Range("E36").Value = 0.001103549751701306
$ws.Range("F36").Value = 0.0003678499172337686
$ws.Range("G36").Value = 0.0003678499172337686
$ws.Range("H36").Value = 0.1286555085525106
$ws.Range("I36").Value = 0.4517196983630679
$ws.Range("J36").Value = 0.3148795291521059
$ws.Range("J37").Value = 0.008276623137759793
$ws.Range("E38").Value = 0.002023174544785727
$ws.Range("F38").Value = 0.001839249586168843
$ws.Range("G38").Value = 0.7544601802464594
$ws.Range("H38").Value = 0.4128195696155968
$ws.Range("I38").Value = 0.0003678499172337686
$ws.Range("J38").Value = 0.001103549751701306
$ws.Range("E39").Value = 0
$ws.Range("H39").Value = 0.003034761817178591
$ws.Range("I39").Value = 0.003310649255103918
$ws.Range("J39").Value = 0.0003678499172337686
$ws.Range("E40").Value = 0.1575317270553614
$ws.Range("F40").Value = 0.1986389553062351
$ws.Range("G40").Value = 0.008644473054993563
$ws.Range("H40").Value = 0.06115504874011403
$ws.Range("I40").Value = 0.05499356262644841
$ws.Range("J40").Value = 0.05122310097480227
$ws.Range("J41").Value = 0.0007356998344675372
$ws.Range("H42").Value = 0.0003678499172337686
$ws.Range("E43").Value = 0.6315983078903807
$ws.Range("F43").Value = 0.5926062166636012
$ws.Range("G43").Value = 0.01121942247562994
$ws.Range("H43").Value = 0.07127092146404268
$ws.Range("I43").Value = 0.2187787382747839
$ws.Range("J43").Value = 0.2451719698363068
$ws.Range("E44").Value = 0.01094353503770462
$ws.Range("F44").Value = 0.00128747471031819
$ws.Range("G44").Value = 0.0003678499172337686
$ws.Range("H44").Value = 0.1116424498804488
$ws.Range("I44").Value = 0.04892403899209122
$ws.Range("J44").Value = 0.008460548096376679
$ws.Range("I45").Value = 0.0001839249586168843
$ws.Range("H46").Value = 0.00009196247930844215
$ws.Range("I46").Value = 0.00009196247930844215
$ws.Range("E47").Value = 0
$ws.Range("H47").Value = 0.00009196247930844215
$ws.Range("E48").Value = 0.02823248114769174
$ws.Range("F48").Value = 0.01563362148243517
$ws.Range("G48").Value = 0.0001839249586168843
$ws.Range("H48").Value = 0.01140334743424683
$ws.Range("I48").Value = 0.06694868493654589
$ws.Range("J48").Value = 0.001103549751701306
$ws.Range("H49").Value = 0
